$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1034.2593
$ws.Range("I28").Value = 222.72223
$ws.Range("J28").Value = 2657.3333
$ws.Range("K28").Value = 222.72223
$ws.Range("L28").Value = 2657.3333
$ws.Range("M28").Value = 262.27777
$ws.Range("N28").Value = -3627.3333

$ws.Range("H137").Value = 1306.6072
$ws.Range("I137").Value = 778.4857
$ws.Range("K137").Value = 2335.4571
$ws.Range("M137").Value = 214.5429000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18432.926
$ws.Range("I32").Value = 20115.365
$ws.Range("J32").Value = 12544.389
$ws.Range("K32").Value = 20115.365
$ws.Range("L32").Value = 12544.389
$ws.Range("M32").Value = -19828.365
$ws.Range("N32").Value = -13118.389

$ws.Range("H45").Value = 47620692
$ws.Range("I45").Value = 111112280
$ws.Range("K45").Value = 111112280
$ws.Range("M45").Value = -111111903

$ws.Range("H61").Value = 1485.2094
$ws.Range("I61").Value = 1121.4546
$ws.Range("J61").Value = 2685.6
$ws.Range("K61").Value = 1121.4546
$ws.Range("L61").Value = 2685.6
$ws.Range("M61").Value = -909.4546
$ws.Range("N61").Value = -3109.6

$ws.Range("H63").Value = 2501949.8
$ws.Range("I63").Value = 3334999.8
$ws.Range("K63").Value = 3334999.8
$ws.Range("M63").Value = -3334313.8

$ws.Range("H66").Value = 2501949.8
$ws.Range("I66").Value = 3334999.8
$ws.Range("K66").Value = 16674999
$ws.Range("M66").Value = -16671567

$ws.Range("H74").Value = 516.5
$ws.Range("I74").Value = 333.25
$ws.Range("J74").Value = 883
$ws.Range("K74").Value = 333.25
$ws.Range("L74").Value = 883
$ws.Range("M74").Value = 540.75
$ws.Range("N74").Value = -2631

$ws.Range("H77").Value = 516.5
$ws.Range("I77").Value = 333.25
$ws.Range("J77").Value = 883
$ws.Range("K77").Value = 1666.25
$ws.Range("L77").Value = 4415
$ws.Range("M77").Value = 2701.75
$ws.Range("N77").Value = -13151

$ws.Range("H122").Value = 2023.6666
$ws.Range("I122").Value = 1809.3334
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 5428.0002
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -2978.0002
$ws.Range("N122").Value = -12900.0001

$ws.Range("H132").Value = 4200.902
$ws.Range("I132").Value = 4460.436
$ws.Range("J132").Value = 3357.4167
$ws.Range("K132").Value = 13381.308
$ws.Range("L132").Value = 10072.2501
$ws.Range("M132").Value = -10851.308
$ws.Range("N132").Value = -15132.2501

$ws.Range("H136").Value = 1485.2094
$ws.Range("I136").Value = 1121.4546
$ws.Range("J136").Value = 2685.6
$ws.Range("K136").Value = 3364.3638
$ws.Range("L136").Value = 8056.799999999999
$ws.Range("M136").Value = -814.3638000000001
$ws.Range("N136").Value = -13156.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4408.913
$ws.Range("I20").Value = 5646
$ws.Range("J20").Value = 3274.9167
$ws.Range("K20").Value = 5646
$ws.Range("L20").Value = 3274.9167
$ws.Range("M20").Value = -5399
$ws.Range("N20").Value = -3768.9167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3790691.8
$ws.Range("I31").Value = 2205.3142
$ws.Range("K31").Value = 2205.3142
$ws.Range("M31").Value = -1910.3142

$ws.Range("H34").Value = 3790691.8
$ws.Range("I34").Value = 2205.3142
$ws.Range("K34").Value = 2205.3142
$ws.Range("M34").Value = -2003.3142

$ws.Range("H58").Value = 1490.8462
$ws.Range("I58").Value = 1297.8889
$ws.Range("K58").Value = 1297.8889
$ws.Range("M58").Value = -1094.8889

$ws.Range("H86").Value = 166669000
$ws.Range("I86").Value = 250001810
$ws.Range("J86").Value = 3360.5
$ws.Range("K86").Value = 250001810
$ws.Range("L86").Value = 3360.5
$ws.Range("M86").Value = -250000687
$ws.Range("N86").Value = -5606.5

$ws.Range("H89").Value = 166669000
$ws.Range("I89").Value = 250001810
$ws.Range("J89").Value = 3360.5
$ws.Range("K89").Value = 1250009050
$ws.Range("L89").Value = 16802.5
$ws.Range("M89").Value = -1250003434
$ws.Range("N89").Value = -28034.5

$ws.Range("H94").Value = 58824588
$ws.Range("I94").Value = 500000400
$ws.Range("J94").Value = 1146.8667
$ws.Range("K94").Value = 500000400
$ws.Range("L94").Value = 1146.8667
$ws.Range("M94").Value = -499999949
$ws.Range("N94").Value = -2048.8667

$ws.Range("H107").Value = 584.2963
$ws.Range("I107").Value = 516.6389
$ws.Range("K107").Value = 516.6389
$ws.Range("M107").Value = 1403.3611

$ws.Range("H132").Value = 2978463.2
$ws.Range("I132").Value = 1631.9
$ws.Range("J132").Value = 10420542
$ws.Range("K132").Value = 4895.700000000001
$ws.Range("L132").Value = 31261626
$ws.Range("M132").Value = -2365.700000000001
$ws.Range("N132").Value = -31266686

$ws.Range("H134").Value = 1149.5834
$ws.Range("I134").Value = 1117.7273
$ws.Range("K134").Value = 3353.1819
$ws.Range("M134").Value = -818.1819

$ws.Range("H136").Value = 1490.8462
$ws.Range("I136").Value = 1297.8889
$ws.Range("K136").Value = 3893.6667
$ws.Range("M136").Value = -1343.6667

$ws.Range("H140").Value = 51666.668
$ws.Range("J140").Value = 35000
$ws.Range("L140").Value = 35000
$ws.Range("N140").Value = -45360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2160
$ws.Range("I81").Value = 320
$ws.Range("K81").Value = 960
$ws.Range("M81").Value = 163

$ws.Range("H84").Value = 2160
$ws.Range("I84").Value = 320
$ws.Range("K84").Value = 2880
$ws.Range("M84").Value = 2736

$ws.Range("H131").Value = 764.15
$ws.Range("J131").Value = 778.07294
$ws.Range("L131").Value = 2334.21882
$ws.Range("N131").Value = -12414.21882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3378.3333
$ws.Range("I80").Value = 3321
$ws.Range("K80").Value = 3321
$ws.Range("M80").Value = -2323

$ws.Range("H83").Value = 3378.3333
$ws.Range("I83").Value = 3321
$ws.Range("K83").Value = 16605
$ws.Range("M83").Value = -11613

$ws.Range("H102").Value = 1836
$ws.Range("I102").Value = 2006
$ws.Range("K102").Value = 2006
$ws.Range("M102").Value = -384

$ws.Range("H113").Value = 16667472
$ws.Range("I113").Value = 125000150
$ws.Range("J113").Value = 904.7692
$ws.Range("K113").Value = 125000150
$ws.Range("L113").Value = 904.7692
$ws.Range("M113").Value = -124997980
$ws.Range("N113").Value = -5244.7692

$ws.Range("H122").Value = 32260350
$ws.Range("I122").Value = 62502036
$ws.Range("K122").Value = 187506108
$ws.Range("M122").Value = -187503658

$ws.Range("H126").Value = 1515.0938
$ws.Range("I126").Value = 1737.7778
$ws.Range("J126").Value = 1228.7858
$ws.Range("K126").Value = 5213.3334
$ws.Range("L126").Value = 3686.3574
$ws.Range("M126").Value = -2743.3334
$ws.Range("N126").Value = -8626.357400000001

$ws.Range("H132").Value = 30394.473
$ws.Range("I132").Value = 43967.293
$ws.Range("J132").Value = 3248.8333
$ws.Range("K132").Value = 131901.879
$ws.Range("L132").Value = 9746.499899999999
$ws.Range("M132").Value = -129371.879
$ws.Range("N132").Value = -14806.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26317620
$ws.Range("I7").Value = 1434.2727
$ws.Range("K7").Value = 1434.2727
$ws.Range("M7").Value = -1322.2727

$ws.Range("H16").Value = 1400
$ws.Range("I16").Value = 1400
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1400
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1230
$ws.Range("N16").ClearContents()

$ws.Range("H100").Value = 47620844
$ws.Range("I100").Value = 83334690
$ws.Range("J100").Value = 2383.3333
$ws.Range("K100").Value = 83334690
$ws.Range("L100").Value = 2383.3333
$ws.Range("M100").Value = -83334149
$ws.Range("N100").Value = -3465.3333

$ws.Range("H126").Value = 26317620
$ws.Range("I126").Value = 1434.2727
$ws.Range("K126").Value = 4302.8181
$ws.Range("M126").Value = -1832.8181

$ws.Range("H132").Value = 2402.6333
$ws.Range("I132").Value = 2121
$ws.Range("J132").Value = 2649.0625
$ws.Range("K132").Value = 6363
$ws.Range("L132").Value = 7947.1875
$ws.Range("M132").Value = -3833
$ws.Range("N132").Value = -13007.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 603.4286
$ws.Range("I107").Value = 587.375
$ws.Range("J107").Value = 624.8333
$ws.Range("K107").Value = 1762.125
$ws.Range("L107").Value = 1874.4999
$ws.Range("M107").Value = 157.875
$ws.Range("N107").Value = -5714.4999

$ws.Range("H126").Value = 83336410
$ws.Range("I126").Value = 90911990
$ws.Range("K126").Value = 272735970
$ws.Range("M126").Value = -272733500

$ws.Range("H136").Value = 1076.2075
$ws.Range("I136").Value = 1083.6216
$ws.Range("K136").Value = 3250.8648
$ws.Range("M136").Value = -700.8647999999998
